$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.833.54'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.626.61'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.519'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.47'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('D12').Value = '1.859.50'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').Value = '1.638.64'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.27'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.06%  '
$ws.Range('D16').Value = '29.859.79'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.110'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('E31').Value = '  +2.53%  '
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('D34').Value = '1.425.37'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('E36').Value = '  -2.72%  '
$ws.Range('E37').Value = '  -4.85%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '74.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.22%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.555'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0498'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').Value = '1.766.33'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '90.73'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.51%  '
$ws.Range('E51').Value = '  +7.51%  '
